# Rewrite the "KEY ACHIEVEMENTS AND IMPACT" bullet list to use concise,
# impact-focused accomplishment statements instead of the longer job-duty
# style bullets (and drop the list from 6 bullets down to 4).
#
# We locate the section by its Heading2 ("KEY ACHIEVEMENTS AND IMPACT"),
# find the first bullet paragraph after the "Impact" sub-heading, and the
# last of the six existing bullet paragraphs, then rewrite that whole
# paragraph range in one shot so Word collapses/expands the paragraph
# count as needed.

$d = $word.ActiveDocument

$bullet = [string][char]0x2022

# --- Locate "KEY ACHIEVEMENTS AND IMPACT" heading paragraph ---
$sectionIndex = -1
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text
    if ($t -match "KEY ACHIEVEMENTS AND IMPACT") {
        $sectionIndex = $i
        break
    }
}
if ($sectionIndex -eq -1) {
    throw "Could not find 'KEY ACHIEVEMENTS AND IMPACT' heading"
}

# --- First bullet paragraph is the first paragraph after the heading
#     (and the "Impact" sub-heading, if present) whose text starts with
#     the bullet character. ---
$firstBulletIndex = -1
for ($i = $sectionIndex + 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($t.StartsWith($bullet)) {
        $firstBulletIndex = $i
        break
    }
}
if ($firstBulletIndex -eq -1) {
    throw "Could not find first bullet under 'KEY ACHIEVEMENTS AND IMPACT'"
}

# --- Last bullet paragraph is the last contiguous paragraph starting
#     with the bullet character (stop at the next non-bullet paragraph,
#     e.g. the following "TECHNICAL SKILLS" heading). ---
$lastBulletIndex = $firstBulletIndex
for ($i = $firstBulletIndex + 1; $i -le $d.Paragraphs.Count; $i++) {
    $t = $d.Paragraphs.Item($i).Range.Text.TrimEnd([char]13, [char]7)
    if ($t.StartsWith($bullet)) {
        $lastBulletIndex = $i
    } else {
        break
    }
}

# --- Build the replacement bullet text (single-line, bullet separated,
#     impact-focused accomplishment statements). ---
$newBullets = @(
    ($bullet + " Platform impact: Built redistricting system serving 12,847 analysts across 89 organizations"),
    ($bullet + " Real-time collaboration at national scale"),
    ($bullet + " Revenue generation: Delivered `$4.9M additional revenue through optimization"),
    ($bullet + " 23% conversion rate improvement")
)

$startRange = $d.Paragraphs.Item($firstBulletIndex).Range.Start
$endRange = $d.Paragraphs.Item($lastBulletIndex).Range.End
$targetRange = $d.Range($startRange, $endRange)

$targetRange.Text = [string]::Join("`r", $newBullets) + "`r"
